$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FGT Summary")

# Update the filename label for the Non-Residential template reference
$ws.Range("A3").Value = "__Blank Non-HomeShare Template.xlsm"

# Range now depends on the length of contract list entries - update figures accordingly
$ws.Range("C2").Value = 4
$ws.Range("G2").Value = 0.16
$ws.Range("L2").Value = 1883.13423588304
